# Refresh the cryptos price sheet ("Updated symbol list" GitHub Actions run).
# Prices are stored as text in column D (General-formatted cells), so a
# leading apostrophe keeps them text-typed and preserves exact formatting
# (e.g. trailing zeros like "0.8140") instead of Excel auto-coercing them
# to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.65"
$ws.Range("D4").Value = "'5.398"
$ws.Range("D5").Value = "'0.06043"
$ws.Range("D6").Value = "'3.397"
$ws.Range("D7").Value = "'0.8140"
$ws.Range("D8").Value = "'0.9248"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01122"
$ws.Range("E9").Value = "8OneONEBestin24h"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1441"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07476"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03369"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03048"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09396"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'4.003"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001597"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04817"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("D18").Value = "'0.005450"
$ws.Range("D19").Value = "'0.004158"
$ws.Range("D20").Value = "'0.0009891"
$ws.Range("D21").Value = "'0.00008802"
$ws.Range("D23").Value = "'6.429"
$ws.Range("D27").Value = "'0.0002901"
$ws.Range("D41").Value = "'0.006402"
$ws.Range("D43").Value = "'0.002901"
$ws.Range("D44").Value = "'0.006380"
$ws.Range("D45").Value = "'0.00005239"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").Value = "'0.002527"
